# Fruta / hortaliza, semanal
# Insert a new weekly price-report row into the Betarraga (Vega Modelo de Temuco) sheet.
# The new record is inserted at row 239, pushing the existing rows 239-322 down to 240-323.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 239 (shifts rows 239..322 -> 240..323)
$ws.Rows(239).Insert()

# Populate the newly inserted row 239 with the new weekly record.
$ws.Range("A239").Value = 10
$ws.Range("B239").Value = "Vega Modelo de Temuco"
$ws.Range("C239").Value = "La Araucanía"
$ws.Range("D239").Value = 44559
$ws.Range("E239").Value = 9
$ws.Range("F239").Value = 100114014
$ws.Range("G239").Value = "Betarraga"
$ws.Range("H239").Value = "Sin especificar"
$ws.Range("I239").Value = "Primera"
$ws.Range("J239").Value = 125
$ws.Range("K239").Value = 600
$ws.Range("L239").Value = 600
$ws.Range("M239").Value = 600
$ws.Range("N239").Value = "$/paquete 5 unidades"
$ws.Range("O239").Value = "Región del Maule"
$ws.Range("P239").Value = 120
$ws.Range("Q239").Value = 5
$ws.Range("R239").Value = "Hortaliza"
